$d = $word.ActiveDocument

# 1. Rename "Summary" heading to "Inverse Relation Among Partial Derivatives"
$d.Content.Find.Execute("Summary", $true, $false, $false, $false, $false, `
                         $true, 1, $false, `
                         "Inverse Relation Among Partial Derivatives", 2)

# 2. Locate that paragraph (it keeps the ListParagraph / ilvl=1 / numId=2 / underline
#    formatting that "Summary" originally had).
$heading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Inverse Relation Among Partial Derivatives") {
        $heading = $cand
    }
}

# 3. Insert the new sub-bullet paragraphs right after it, in order. Each new
#    paragraph inherits the formatting (pPr/rPr) of the one it follows, which
#    matches the target markup exactly.
$titles = @(
    "Ratio Relation for Partials",
    "Relation for Partials with Three Variables",
    "Chain Rule",
    "Exact Equation Criteria",
    "Legendre Transformations",
    ("Stirling" + [char]0x2019 + "s Approximation"),
    "Other Approximations"
)

$prev = $heading
foreach ($title in $titles) {
    $prev.Range.InsertParagraphAfter()
    $next = $prev.Next()
    $next.Range.Text = $title
    $prev = $next
}

# 4. Move the "_GoBack" bookmark from the end of the (renamed) heading paragraph
#    to the start of the "Relevant Examples" paragraph.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.TrimEnd([char]13, [char]7) -eq "Relevant Examples") {
        $target = $cand
    }
}
$r = $target.Range.Duplicate()
$r.Collapse(1)
$d.Bookmarks.Add("_GoBack", $r)
